# Weekly crime-data refresh for CompStat_1 (108th Precinct), week of
# 11/27/2023 - 12/3/2023 (Volume 30, Number 48).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: issue number and the reporting week's date range.
# ---------------------------------------------------------------------
$ws.Range("A8").Value  = "Volume 30   Number  48"
$ws.Range("C9").Value  = "Report Covering the Week  11/27/2023  Through  12/3/2023"

# ---------------------------------------------------------------------
# Crime Complaints table (rows 14-30, columns C:N)
# ---------------------------------------------------------------------

# Row 15 - Rape
$ws.Range("N15").Value = 100

# Row 16 - Robbery
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = 18.75
$ws.Range("I16").Value = 241
$ws.Range("J16").Value = 181
$ws.Range("K16").Value = 33.149171270718
$ws.Range("L16").Value = 109.565217391304
$ws.Range("M16").Value = 35.393258426966
$ws.Range("N16").Value = -74.388947927736

# Row 17 - Fel. Assault
$ws.Range("C17").Value = 6
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = 7.142857142857
$ws.Range("I17").Value = 224
$ws.Range("J17").Value = 206
$ws.Range("K17").Value = 8.737864077669
$ws.Range("L17").Value = 16.062176165803
$ws.Range("M17").Value = 85.123966942148
$ws.Range("N17").Value = -23.287671232876

# Row 18 - Burglary
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -14.285714285714
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = -16.666666666666
$ws.Range("I18").Value = 200
$ws.Range("J18").Value = 159
$ws.Range("K18").Value = 25.786163522012
$ws.Range("L18").Value = 16.279069767441
$ws.Range("M18").Value = -14.893617021276
$ws.Range("N18").Value = -86.139986139986

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 19
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = 35.714285714285
$ws.Range("F19").Value = 55
$ws.Range("G19").Value = 41
$ws.Range("H19").Value = 34.146341463414
$ws.Range("I19").Value = 691
$ws.Range("J19").Value = 638
$ws.Range("K19").Value = 8.307210031347
$ws.Range("L19").Value = 58.486238532110
$ws.Range("M19").Value = 59.584295612009
$ws.Range("N19").Value = -19.181286549707

# Row 20 - G.L.A.
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -25
$ws.Range("F20").Value = 20
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = 33.333333333333
$ws.Range("I20").Value = 286
$ws.Range("J20").Value = 216
$ws.Range("K20").Value = 32.407407407407
$ws.Range("L20").Value = 68.235294117647
$ws.Range("M20").Value = 50.526315789473
$ws.Range("N20").Value = -85.280494081317

# Row 21 - TOTAL
$ws.Range("C21").Value = 37
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = 19.354838709677
$ws.Range("F21").Value = 126
$ws.Range("G21").Value = 105
$ws.Range("H21").Value = 20
$ws.Range("I21").Value = 1666
$ws.Range("J21").Value = 1420
$ws.Range("K21").Value = 17.323943661971
$ws.Range("L21").Value = 51.454545454545
$ws.Range("M21").Value = 41.90800681431
$ws.Range("N21").Value = -69.676010192937

# Row 22 - Transit (D22/E22 become "no prior-year data" placeholders)
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 7
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 75
$ws.Range("I22").Value = 81
$ws.Range("K22").Value = 8
$ws.Range("L22").Value = 102.5
$ws.Range("M22").Value = 76.086956521739

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 50
$ws.Range("D24").Value = 40
$ws.Range("E24").Value = 25
$ws.Range("F24").Value = 191
$ws.Range("G24").Value = 149
$ws.Range("H24").Value = 28.187919463087
$ws.Range("I24").Value = 1941
$ws.Range("J24").Value = 1423
$ws.Range("K24").Value = 36.401967673928
$ws.Range("L24").Value = 61.75
$ws.Range("M24").Value = 122.0823798627

# Row 25 - Misd. Assault
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 14.285714285714
$ws.Range("F25").Value = 42
$ws.Range("G25").Value = 42
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 483
$ws.Range("J25").Value = 520
$ws.Range("K25").Value = -7.115384615384
$ws.Range("L25").Value = 10.273972602739
$ws.Range("M25").Value = 3.426124197002

# Row 27 - Other Sex Crimes
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -66.666666666666
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 9
$ws.Range("H27").Value = -33.333333333333
$ws.Range("I27").Value = 78
$ws.Range("J27").Value = 87
$ws.Range("K27").Value = -10.344827586206
$ws.Range("L27").Value = 30

# Row 30 - Hate Crimes
$ws.Range("F30").Value = 1
$ws.Range("I30").Value = 11
$ws.Range("K30").Value = 22.222222222222
$ws.Range("L30").Value = -8.333333333333

# ---------------------------------------------------------------------
# Rows 22/28/29: D and E drop from numeric counts/pct-chg to the
# "no activity" placeholders ("0" / "***.*") used elsewhere in the sheet.
# These look like plain numbers, so we briefly force Text formatting to
# stop Excel from re-interpreting "0" as numeric zero, then restore the
# normal (General) formatting used by the surrounding label cells.
# ---------------------------------------------------------------------
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"

$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E29").PasteSpecial(-4122)

$ws.Range("E22").Value = "***.*"
$ws.Range("E28").Value = "***.*"
$ws.Range("E29").Value = "***.*"

Write-Host "Weekly crime data refresh applied."
